$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row (2-244).
for ($row = 2; $row -le 244; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
